$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for data rows 2-99
# from serial date 45174 (2023-09-05) to 45175 (2023-09-06).
$ws.Range("C2:C99").Value = 45175
